$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$ws.Rows.Item(20).Delete()
